$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value()
    $eVal = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($r, 4).Value = $eVal
    $ws.Cells.Item($r, 5).Value = $dVal
}
